# Apply updated cryptocurrency price/volume data to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.673.26"
$ws.Range("E2").Value = "  -6.32%  "
$ws.Range("D3").Value = "2.897.35"
$ws.Range("E3").Value = "  -4.49%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "549.17"
$ws.Range("E5").Value = "  -4.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "121.59"
$ws.Range("E6").Value = "  -6.09%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").Value = "2.882.81"
$ws.Range("E8").Value = "  -4.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.493"
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("E10").Value = "  -8.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.73"
$ws.Range("E11").Value = "  -9.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000210"
$ws.Range("E13").Value = "  -8.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.59"
$ws.Range("E14").Value = "  -5.49%  "
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("D16").Value = "3.373.64"
$ws.Range("E16").Value = "  -4.53%  "
$ws.Range("D17").Value = "2.890.19"
$ws.Range("E17").Value = "  -4.71%  "
$ws.Range("D18").Value = "57.603.65"
$ws.Range("E18").Value = "  -6.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.48"
$ws.Range("E19").Value = "  +2.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "406.85"
$ws.Range("E20").Value = "  -8.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.81"
$ws.Range("E21").Value = "  -4.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.653"
$ws.Range("E22").Value = "  -2.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.54"
$ws.Range("E24").Value = "  -1.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "76.83"
$ws.Range("E25").Value = "  -4.39%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.45"
$ws.Range("E28").Value = "  -3.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.16"
$ws.Range("E30").Value = "  -4.24%  "
$ws.Range("E31").Value = "  -5.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "24.62"
$ws.Range("E32").Value = "  -4.23%  "
$ws.Range("E33").Value = "  -0.85%  "
$ws.Range("E34").Value = "  -12.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.896"
$ws.Range("E35").Value = "  -7.49%  "
$ws.Range("E36").Value = "  -5.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "48.37"
$ws.Range("E37").Value = "  -3.87%  "
$ws.Range("E38").Value = "  +6.68%  "
$ws.Range("E39").Value = "  -11.59%  "
$ws.Range("E40").Value = "  -7.58%  "
$ws.Range("E41").Value = "  -4.53%  "
$ws.Range("D42").Value = "2.607.53"
$ws.Range("E42").Value = "  -2.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "359.02"
$ws.Range("E43").Value = "  -4.43%  "
$ws.Range("E44").Value = "  -6.99%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "117.51"
$ws.Range("E46").Value = "  -4.89%  "
$ws.Range("E48").Value = "  -1.30%  "
$ws.Range("E49").Value = "  -3.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.38"
$ws.Range("E50").Value = "  -5.08%  "
$ws.Range("E51").Value = "  -5.21%  "
